$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells below store plain numbers/percentages as literal TEXT (inline
# strings) in the source workbook. Assigning a numeric-looking string to a
# General-formatted cell makes Excel auto-convert it to a real number, so we
# first mark each target cell as Text ("@") and then write the new literal
# value, which keeps it stored as text exactly like the original cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "288.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.71%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.30%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.268"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.40%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06991"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.36%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.458"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.56%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.556"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.12%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.390"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.86%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9041"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.72%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1603"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.34%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07586"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "14.48%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07724"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.95%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02913"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.07%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09032"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.46%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001591"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.14%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006520"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.13%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006422"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.14%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.491"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.13%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.232"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.63%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.01%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1336"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.04%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.016"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.77%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1596"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.05%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04526"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.78%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.45%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004147"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.77%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001169"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.22%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001668"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "3.37%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04362"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006930"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.95%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1250"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.42%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002068"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.67%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01163"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.97%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005824"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.17%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.26%"
